# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the crypto symbol
# list with updated quote data. Values are written with a leading
# apostrophe so Excel stores them as literal text (preserving formatting
# like trailing zeros and the "%" suffix) instead of re-parsing them as
# numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.30"
$ws.Range("E2").Value = "'-0.58%"
$ws.Range("D3").Value = "'27.17"
$ws.Range("E3").Value = "'2.88%"
$ws.Range("D4").Value = "'5.091"
$ws.Range("E4").Value = "'0.05%"
$ws.Range("D5").Value = "'0.05706"
$ws.Range("E5").Value = "'1.89%"
$ws.Range("D6").Value = "'6.497"
$ws.Range("E6").Value = "'0.20%"
$ws.Range("D7").Value = "'0.8190"
$ws.Range("D8").Value = "'0.8642"
$ws.Range("E8").Value = "'2.14%"
$ws.Range("E9").Value = "'-0.59%"
$ws.Range("D10").Value = "'0.06929"
$ws.Range("E10").Value = "'-1.07%"
$ws.Range("D11").Value = "'0.02830"
$ws.Range("E11").Value = "'-0.60%"
$ws.Range("D12").Value = "'0.09396"
$ws.Range("E12").Value = "'0.09%"
$ws.Range("D13").Value = "'0.001521"
$ws.Range("E13").Value = "'0.27%"
$ws.Range("D14").Value = "'0.04059"
$ws.Range("E14").Value = "'-12.93%"
$ws.Range("E15").Value = "'0.44%"
$ws.Range("D16").Value = "'0.006173"
$ws.Range("E16").Value = "'0.52%"
$ws.Range("D17").Value = "'3.508"
$ws.Range("E17").Value = "'-2.64%"
$ws.Range("D18").Value = "'3.009"
$ws.Range("E18").Value = "'-0.16%"
$ws.Range("D19").Value = "'2.228"
$ws.Range("E19").Value = "'8.41%"
$ws.Range("D20").Value = "'0.3165"
$ws.Range("E20").Value = "'-1.29%"
$ws.Range("D21").Value = "'0.03198"
$ws.Range("E21").Value = "'0.54%"
$ws.Range("E22").Value = "'0.50%"
$ws.Range("D23").Value = "'3.561"
$ws.Range("E23").Value = "'-4.88%"
$ws.Range("E25").Value = "'-2.50%"
$ws.Range("D26").Value = "'0.004482"
$ws.Range("E26").Value = "'-2.20%"
$ws.Range("D27").Value = "'0.00009891"
$ws.Range("E27").Value = "'3.03%"
$ws.Range("E28").Value = "'3.61%"
$ws.Range("D40").Value = "'0.03727"
$ws.Range("E40").Value = "'1.62%"
$ws.Range("D41").Value = "'0.005904"
$ws.Range("E41").Value = "'-4.92%"
$ws.Range("D42").Value = "'0.1060"
$ws.Range("E42").Value = "'0.32%"
$ws.Range("D43").Value = "'0.002488"
$ws.Range("E43").Value = "'-0.50%"
$ws.Range("D44").Value = "'0.009388"
$ws.Range("E44").Value = "'1.54%"
$ws.Range("D45").Value = "'0.00005165"
$ws.Range("E45").Value = "'-2.31%"
$ws.Range("E46").Value = "'-0.07%"
$ws.Range("E47").Value = "'7.02%"
$ws.Range("D48").Value = "'0.002517"
$ws.Range("E48").Value = "'-4.53%"
$ws.Range("E49").Value = "'-0.07%"
$ws.Range("E50").Value = "'-0.07%"
